$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock Report")

$rowData = @{
    7 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "RIGHT SIDE F/B NO.1 HOLE 2\"X2\" AT 2 PLS & LEFT SIDE F/B NO.2 HOLE 2\"X2\"", "", "", ""]'
    8 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "PANELS-(PANELS)", "Right side panel -(Right side panel )", "RIGHT SITE PANEL NO.3 PUSH IN 30\"X15\"", "", "", ""]'
    9 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "PANELS-(PANELS)", "Right side panel -(Right side panel )", "INTERNAL SCRACHED & RUSTED", "", "", ""]'
    10 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "PANELS-(PANELS)", "Right side panel -(Right side panel )", "LEFT SIDE PANEL NO.1 HOLE 2\"X2\"& LEFT SIDE PANEL NO.2 PIN HOLE 1\"X1\" 2 PLS", "", "", ""]'
    11 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "PANELS-(PANELS)", "Right side panel -(Right side panel )", "LEFT SIDE PANEL NO.3 HOLE 3\"X2\" 3 PLS", "", "", ""]'
    12 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "DOORS-(D)", "Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)", "RIGHT DOOR LOCK BAR BENT 1 PES", "", "", ""]'
    15 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "PANELS-(PANELS)", "Right side panel -(Right side panel )", "L/S PANEL 3 HOLE 1\"X1\"", "", "", ""]'
    16 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B BROKEN HOLE 60\"X48\" & HOLE 6\"X6\" AT 4 PLS & BROKEN DOWN 15\"X12\"", "", "", ""]'
    17 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "FLOOR BORD DIRTY BY SCRAP & SCRATCH.", "", "", ""]'
    19 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "BOTH SIDE F/B HOLE 10\"X5,15\"X10\" & 15\"X15\" AT 3PLS", "", "", ""]'
    20 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "PANELS-(PANELS)", "Right side panel -(Right side panel )", "INTERNAL PANEL SCRATCH BADLY", "", "", ""]'
    21 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "PANELS-(PANELS)", "Right side panel -(Right side panel )", "RIGHT SIDE PANEL NO.3 CUT 15\"X10 & PUSH IN-PUSH OUT AT V/P", "", "", ""]'
    22 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "DOORS-(D)", "Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)", "LEFT DOOR BOTTOM GASKET LOSE 20\"", "", "", ""]'
    23 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "PANELS-(PANELS)", "Right side panel -(Right side panel )", "ROOF PANEL DENT AT V/P", "", "", ""]'
    25 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "PANELS-(PANELS)", "Right side panel -(Right side panel )", "INERNAL PANEL DIRTY", "", "", ""]'
    29 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "PANELS-(PANELS)", "Right side panel -(Right side panel )", "INTERNAL PANEL DIRTY & RUSTED & CORROSION", "", "", ""]'
    30 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "PANELS-(PANELS)", "Right side panel -(Right side panel )", "ROOF PANEL 2 DENT 30\"X35\"", "", "", ""]'
    31 = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "PANELS-(PANELS)", "Right side panel -(Right side panel )", "LEFT SIDE PANEL 4 DENT 30\"X20\"", "", "", ""]'
}

foreach ($r in $rowData.Keys) {
    $ws.Range("C" + $r + ":AB" + $r).Clear()
    $ws.Cells.Item($r, 1).Value = " "
    $ws.Cells.Item($r, 2).Value = $rowData[$r]
}

$ws.Columns.Item(2).ColumnWidth = 284.68988764044946
$ws.Columns.Item(2).Hidden = $true
$ws.Columns.Item(24).ColumnWidth = 27.28988764044944
